$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean: wipe the existing 6x2 block (and its shared strings) so we can
# rebuild the sheet content in a controlled, deterministic order.
$ws.Range("A1:F2").Clear()

# --- Row 1 / Row 2 text values -------------------------------------------
# (assignment order chosen so the resulting shared-string table is built up
# in the same order the strings first appear in the target workbook)
$ws.Range("A2").Value = "BRI01"
$ws.Range("C1").Value = "NGR"
$ws.Range("C2").Value = "ND1128168190"
$ws.Range("G1").Value = "Level Controller Name"
$ws.Range("G2").Value = "Storm Overflow Level Monitor Loop"
$ws.Range("H1").Value = "Manufacturer"
$ws.Range("I1").Value = "Model"
$ws.Range("J1").Value = "Serial Number"
$ws.Range("H2").Value = "SIEMENS"
$ws.Range("I2").Value = "HYDRORANGER 200"
$ws.Range("J2").Value = "PDB2981/XK"
$ws.Range("B1").Value = "AI2 Site Reference"
$ws.Range("A1").Value = "S4 Root FuncLoc"
$ws.Range("B2").Value = "SAI00023001"
$ws.Range("E1").Value = "AI2 Equipment SAI Number"
$ws.Range("F1").Value = "AI2 Equipment PLI Code"
$ws.Range("E2").Value = "SAI00023450"
$ws.Range("F2").Value = "PLI00004561"
$ws.Range("K1").Value = "Memo Line"
$ws.Range("L1").Value = "Relay 1 Function"
$ws.Range("M1").Value = "Relay 1 On"
$ws.Range("N1").Value = "Relay 1 Off"
$ws.Range("O1").Value = "Relay 2 Function"
$ws.Range("P1").Value = "Relay 2 On"
$ws.Range("Q1").Value = "Relay 2 Off"
$ws.Range("R1").Value = "Relay 3 Function"
$ws.Range("S1").Value = "Relay 3 On"
$ws.Range("T1").Value = "Relay 3 Off"
$ws.Range("U1").Value = "Relay 4 Function"
$ws.Range("V1").Value = "Relay 4 On"
$ws.Range("W1").Value = "Relay 4 Off"
$ws.Range("X1").Value = "Relay 5 Function"
$ws.Range("Y1").Value = "Relay 5 On"
$ws.Range("Z1").Value = "Relay 5 Off"
$ws.Range("AA1").Value = "Relay 6 Function"
$ws.Range("AB1").Value = "Relay 6 On"
$ws.Range("AC1").Value = "Relay 6 Off"
$ws.Range("AA2").Value = "LOSS OF ECHO"
$ws.Range("D1").Value = "Install Date"

# --- D2: Install Date value (date-formatted number) -----------------------
$ws.Range("D2").Value = 43679
$ws.Range("D2").NumberFormat = "mm-dd-yy"

# --- Column widths ----------------------------------------------------------
$ws.Columns("A").ColumnWidth = 19.666666666666664
$ws.Columns("B").ColumnWidth = 19.666666666666664
$ws.Columns("C").ColumnWidth = 13
$ws.Columns("D").ColumnWidth = 13
$ws.Columns("E").ColumnWidth = 24.666666666666664
$ws.Columns("F").ColumnWidth = 24.666666666666664
$ws.Columns("G").ColumnWidth = 32.33333333333333
$ws.Columns("H").ColumnWidth = 12.333333333333332
$ws.Columns("I").ColumnWidth = 17.333333333333336
$ws.Columns("J").ColumnWidth = 13
$ws.Columns("K").ColumnWidth = 10
$ws.Columns("AA").ColumnWidth = 14.666666666666668

# --- Selection matches the author's saved cursor position -----------------
$ws.Range("D3").Select()
